$wb = $excel.ActiveWorkbook

# Remove the stray "Sheet" diagnostic row from optimization_parameters
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()

# Restore the selection left behind on that sheet
$ws.Range("C41").Select()

# Make optimization_diagnostics the active sheet/tab
$ws2 = $wb.Worksheets.Item("optimization_diagnostics")
$ws2.Activate()
